$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (periods reordered newest-first, amounts refreshed) for the
# two workers in the "Estado de Cuenta" table (rows 16-29, cols B:G).
$data = @(
    @{ Row = 16; B = "CC"; C = "73148203";   D = "PEDRO UTRIA MONSALVE";          E = "2001"; F = 70000;  G = 877803 },
    @{ Row = 17; B = "CC"; C = "73148203";   D = "PEDRO UTRIA MONSALVE";          E = "1912"; F = 100000; G = 877803 },
    @{ Row = 18; B = "CC"; C = "73148203";   D = "PEDRO UTRIA MONSALVE";          E = "1911"; F = 100000; G = 877803 },
    @{ Row = 19; B = "CC"; C = "73148203";   D = "PEDRO UTRIA MONSALVE";          E = "1910"; F = 100000; G = 877803 },
    @{ Row = 20; B = "CC"; C = "73148203";   D = "PEDRO UTRIA MONSALVE";          E = "1909"; F = 100000; G = 877803 },
    @{ Row = 21; B = "CC"; C = "73148203";   D = "PEDRO UTRIA MONSALVE";          E = "1908"; F = 100000; G = 877803 },
    @{ Row = 22; B = "CC"; C = "73148203";   D = "PEDRO UTRIA MONSALVE";          E = "1907"; F = 100000; G = 877803 },
    @{ Row = 23; B = "CC"; C = "1152226681"; D = "LUIS MIGUELL UTRIA MORALESS";   E = "2001"; F = 23187;  G = 828116 },
    @{ Row = 24; B = "CC"; C = "1152226681"; D = "LUIS MIGUELL UTRIA MORALESS";   E = "1912"; F = 33125;  G = 828116 },
    @{ Row = 25; B = "CC"; C = "1152226681"; D = "LUIS MIGUELL UTRIA MORALESS";   E = "1911"; F = 33125;  G = 828116 },
    @{ Row = 26; B = "CC"; C = "1152226681"; D = "LUIS MIGUELL UTRIA MORALESS";   E = "1910"; F = 33125;  G = 828116 },
    @{ Row = 27; B = "CC"; C = "1152226681"; D = "LUIS MIGUELL UTRIA MORALESS";   E = "1909"; F = 33125;  G = 828116 },
    @{ Row = 28; B = "CC"; C = "1152226681"; D = "LUIS MIGUELL UTRIA MORALESS";   E = "1908"; F = 33125;  G = 828116 },
    @{ Row = 29; B = "CC"; C = "1152226681"; D = "LUIS MIGUELL UTRIA MORALESS";   E = "1907"; F = 33125;  G = 828116 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

# Column widths widened slightly to keep the bestFit look after the data refresh
# (the new values are a little longer/shorter than the originals).
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
